$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I0 (column I) is 1 for most rows, except the last three (rows 22-24)
$i0Values = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 7, 6, 3)

# IF (column J) mirrors IP (column H) for most rows, except the last three (rows 22-24)
$ifValues = @(3, 5, 5, 5, 6, 6, 4, 5, 6, 1, 6, 6, 6, 5, 5, 5, 5, 5, 5, 4, 9, 8, 4)

for ($idx = 0; $idx -lt 23; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
